$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (102-146) being appended to the master-reg_center_device_h table.
# Columns: A=regcntr_id, B=device_id, C=lang_code("eng"), D=is_active(TRUE),
#          E=cr_by("superadmin"), F=cr_dtimes("now()"), G=eff_dtimes("now()")
$rows = @(
    @(102, 10002, 3000121),
    @(103, 10003, 3000122),
    @(104, 10004, 3000123),
    @(105, 10005, 3000124),
    @(106, 10006, 3000125),
    @(107, 10007, 3000126),
    @(108, 10008, 3000127),
    @(109, 10009, 3000128),
    @(110, 10010, 3000129),
    @(111, 10002, 3000130),
    @(112, 10003, 3000131),
    @(113, 10004, 3000132),
    @(114, 10005, 3000133),
    @(115, 10006, 3000134),
    @(116, 10007, 3000135),
    @(117, 10008, 3000136),
    @(118, 10009, 3000137),
    @(119, 10010, 3000138),
    @(120, 10002, 3000139),
    @(121, 10003, 3000140),
    @(122, 10004, 3000141),
    @(123, 10005, 3000142),
    @(124, 10006, 3000143),
    @(125, 10007, 3000144),
    @(126, 10008, 3000145),
    @(127, 10009, 3000146),
    @(128, 10010, 3000147),
    @(129, 10002, 3000148),
    @(130, 10003, 3000149),
    @(131, 10004, 3000150),
    @(132, 10005, 3000151),
    @(133, 10006, 3000152),
    @(134, 10007, 3000153),
    @(135, 10008, 3000154),
    @(136, 10009, 3000155),
    @(137, 10010, 3000156),
    @(138, 10002, 3000157),
    @(139, 10003, 3000158),
    @(140, 10004, 3000159),
    @(141, 10005, 3000160),
    @(142, 10006, 3000161),
    @(143, 10007, 3000162),
    @(144, 10008, 3000163),
    @(145, 10009, 3000164),
    @(146, 10010, 3000165)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = "eng"
    $ws.Cells.Item($rowNum, 4).Value = $true
    $ws.Cells.Item($rowNum, 5).Value = "superadmin"
    $ws.Cells.Item($rowNum, 6).Value = "now()"
    $ws.Cells.Item($rowNum, 7).Value = "now()"
}

# Update the page setup to portrait, matching the author's print settings.
$ws.PageSetup.Orientation = 1

# Scroll the view down and select the newly-added block, matching the
# author's last on-screen selection before saving.
$win = $wb.Windows.Item(1)
$win.ScrollRow = 128
$win.ScrollColumn = 1
$ws.Range("A102:B146").Select()
